$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000001099251006220214
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 1594453305621061120
$ws.Range("G2").Value = 1594453305621063936

$ws.Range("B3").Value = 0.06328177979961902
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 1.053659104900323
